# Lecture 6 / Multidimensional arrays slide: the "three-dimensional array"
# example currently declares the variable with only two pairs of brackets
# (int[][] arrayInThreeDimension). Fix the code sample so it actually
# declares a 3D array: int[][][] arrayInThreeDimension.

$p = $ppt.ActivePresentation

# Find the "Многомерни масиви" (Multidimensional arrays) slide by locating
# the shape that contains both code samples, rather than hard-coding a
# slide index.
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text.IndexOf("arrayInThreeDimension") -ge 0) {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$full = $tr.Text

# Locate the "[][] " run that immediately precedes the
# "arrayInThreeDimension" identifier (the three-dimensional array sample),
# as opposed to the "arrayInTwoDimension" sample earlier in the same shape.
$marker = "arrayInThreeDimension"
$markerIdx = $full.IndexOf($marker)

$old = "[][] "
$new = "[][][] "

$searchRegion = $full.Substring(0, $markerIdx)
$targetIdx = $searchRegion.LastIndexOf($old)

$chars = $tr.Characters($targetIdx + 1, $old.Length)
$chars.Text = $new
